# Swag Labs TestNG login suite - initial commit
#
# Rebuilds the credentials sheet:
#   - renames the default sheet to "data"
#   - writes a bold, centered "username"/"password" header row
#   - writes the "standard_user"/"secret_sauce" credential row below it
#   - wraps + vertically centers both rows and doubles the row height
#   - leaves the selection on B2, matching the saved workbook state

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename "Sheet1" -> "data"
$ws.Name = "data"

# Header row (row 1): username / password
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"

# Data row (row 2): standard_user / secret_sauce
$ws.Range("A2").Value = "standard_user"
$ws.Range("B2").Value = "secret_sauce"

# Header formatting: wrapped, centered both ways, bold
$header = $ws.Range("A1:B1")
$header.WrapText = $true
$header.VerticalAlignment = -4108     # xlCenter
$header.HorizontalAlignment = -4108   # xlCenter
$header.Font.Bold = $true

# Data row formatting: wrapped, vertically centered
$data = $ws.Range("A2:B2")
$data.WrapText = $true
$data.VerticalAlignment = -4108       # xlCenter

# Both rows now wrap onto two lines - double the row height to match
$ws.Rows.Item(1).RowHeight = 28.8
$ws.Rows.Item(2).RowHeight = 28.8

# Final selection/active cell is B2
$ws.Range("B2").Select()
